# Update view/click counts in column F across the four sheets of the
# "广州-漫展信息" workbook, matching a refreshed data scrape.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 1028
$ws1.Range("F11").Value = 1385
$ws1.Range("F13").Value = 514
$ws1.Range("F15").Value = 1373
$ws1.Range("F16").Value = 821
$ws1.Range("F18").Value = 1425
$ws1.Range("F22").Value = 26
$ws1.Range("F23").Value = 415
$ws1.Range("F24").Value = 33
$ws1.Range("F25").Value = 3589
$ws1.Range("F28").Value = 1584
$ws1.Range("F29").Value = 32

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 40

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 21

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 21
$ws4.Range("F8").Value  = 40
$ws4.Range("F16").Value = 1028
$ws4.Range("F22").Value = 1385
$ws4.Range("F24").Value = 514
$ws4.Range("F26").Value = 1373
$ws4.Range("F27").Value = 821
$ws4.Range("F29").Value = 1425
$ws4.Range("F35").Value = 26
$ws4.Range("F36").Value = 415
$ws4.Range("F37").Value = 33
$ws4.Range("F38").Value = 3589
$ws4.Range("F41").Value = 1584
$ws4.Range("F44").Value = 32
